# Generate Report for Handback
#
# Localization files came back "in sync with en-US" so the status/report
# sheets need to reflect the handback: the Status column moves from
# "Ready for handoff" to "Handed back: in sync with en-US", and each
# language sheet grows two new columns (F = Latest Target File, G = Latest
# Handback File) with hyperlinks to the file that was handed back, plus the
# Latest Handback DateTime (column H) gets stamped with the real time of
# the handback instead of the zero-date placeholder.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

$mdUrlD60 = "https://github.com/OpenLocalizationTest/oltest/blob/cb1a63928fc863f9190ebdce8053f33cdce2a0ae/e2e/d60ccee0-1a26-4daa-ad98-0ae46a7ea6fb.md"
$mdUrlDa9 = "https://github.com/OpenLocalizationTest/oltest/blob/cb1a63928fc863f9190ebdce8053f33cdce2a0ae/e2e/da9d4636-d718-4a62-88ed-1fb2b06e14fb.md"

$mdNameD60 = "d60ccee0-1a26-4daa-ad98-0ae46a7ea6fb.md"
$mdNameDa9 = "da9d4636-d718-4a62-88ed-1fb2b06e14fb.md"

# ── zh-cn ────────────────────────────────────────────────────────────────
$ws = $wb.Worksheets.Item("zh-cn")

$xlfUrlD60zh = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4edfde31e72067446eed6d89ce57dfb0722b4171/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/d60ccee0-1a26-4daa-ad98-0ae46a7ea6fb.5d60ec18ad777d681efeccb6b9f26fcecd9b0184.zh-cn.xlf"
$xlfUrlDa9zh = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/4edfde31e72067446eed6d89ce57dfb0722b4171/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/da9d4636-d718-4a62-88ed-1fb2b06e14fb.154179887020df0f6d9f3e43bb4c4e9d591201a3.zh-cn.xlf"

$xlfNameD60zh = "d60ccee0-1a26-4daa-ad98-0ae46a7ea6fb.5d60ec18ad777d681efeccb6b9f26fcecd9b0184.zh-cn.xlf"
$xlfNameDa9zh = "da9d4636-d718-4a62-88ed-1fb2b06e14fb.154179887020df0f6d9f3e43bb4c4e9d591201a3.zh-cn.xlf"

# Status: "Ready for handoff" -> "Handed back: in sync with en-US"
$ws.Range("C2").Value = $newStatus
$ws.Range("C3").Value = $newStatus

# Latest Handback DateTime: was the zero-date placeholder, now stamped.
$ws.Range("H2").Value = "2016-03-19 08:51:05"
$ws.Range("H3").Value = "2016-03-19 08:51:05"

# New columns F (Latest Target File) / G (Latest Handback File), with
# hyperlinks matching the Source File / Latest Handoff File columns.
$ws.Hyperlinks.Add($ws.Range("F2"), $mdUrlD60, "", "", $mdNameD60)
$ws.Hyperlinks.Add($ws.Range("G2"), $xlfUrlD60zh, "", "", $xlfNameD60zh)
$ws.Hyperlinks.Add($ws.Range("F3"), $mdUrlDa9, "", "", $mdNameDa9)
$ws.Hyperlinks.Add($ws.Range("G3"), $xlfUrlDa9zh, "", "", $xlfNameDa9zh)

# ── de-de ────────────────────────────────────────────────────────────────
$ws = $wb.Worksheets.Item("de-de")

$xlfUrlD60de = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d9480dac98c276b8ec77f13fd79e995312b6364f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/d60ccee0-1a26-4daa-ad98-0ae46a7ea6fb.5d60ec18ad777d681efeccb6b9f26fcecd9b0184.de-de.xlf"
$xlfUrlDa9de = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d9480dac98c276b8ec77f13fd79e995312b6364f/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/da9d4636-d718-4a62-88ed-1fb2b06e14fb.154179887020df0f6d9f3e43bb4c4e9d591201a3.de-de.xlf"

$xlfNameD60de = "d60ccee0-1a26-4daa-ad98-0ae46a7ea6fb.5d60ec18ad777d681efeccb6b9f26fcecd9b0184.de-de.xlf"
$xlfNameDa9de = "da9d4636-d718-4a62-88ed-1fb2b06e14fb.154179887020df0f6d9f3e43bb4c4e9d591201a3.de-de.xlf"

# Status: "Ready for handoff" -> "Handed back: in sync with en-US"
$ws.Range("C2").Value = $newStatus
$ws.Range("C3").Value = $newStatus

# Latest Handback DateTime: was the zero-date placeholder, now stamped.
$ws.Range("H2").Value = "2016-03-19 08:51:10"
$ws.Range("H3").Value = "2016-03-19 08:51:10"

# New columns F (Latest Target File) / G (Latest Handback File), with
# hyperlinks matching the Source File / Latest Handoff File columns.
$ws.Hyperlinks.Add($ws.Range("F2"), $mdUrlD60, "", "", $mdNameD60)
$ws.Hyperlinks.Add($ws.Range("G2"), $xlfUrlD60de, "", "", $xlfNameD60de)
$ws.Hyperlinks.Add($ws.Range("F3"), $mdUrlDa9, "", "", $mdNameDa9)
$ws.Hyperlinks.Add($ws.Range("G3"), $xlfUrlDa9de, "", "", $xlfNameDa9de)

Write-Host "Handback report generated."
